# Update the "Metadata" sheet with the new URL, Version, Date and Publisher.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-snapshot-provider-zip-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Clear the stray ele-1/ext-1 constraint text that used to sit on the
# top-level "Extension" row of the "Elements" sheet (it now only applies
# to the "Extension.extension" row).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
